# Refresh cryptos list snapshot (prices / 1h volume %) - GitHub Actions data sync
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Writes $Val into $CellRef while forcing a text cell (avoids Excel
    # auto-converting numeric-looking strings like "0.0800" into numbers
    # and stripping significant trailing zeros), then restores the default
    # "Normal" style so no stray per-cell formatting is introduced.
    param($Sheet, $CellRef, $Val)
    $c = $Sheet.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '37.151.46'
Set-TextValue $ws 'E2' '  -0.20%  '
Set-TextValue $ws 'D3' '2.074.34'
Set-TextValue $ws 'E3' '  -1.11%  '
Set-TextValue $ws 'E4' '  +0.07%  '
Set-TextValue $ws 'D5' '253.08'
Set-TextValue $ws 'E5' '  +0.78%  '
Set-TextValue $ws 'E6' '  +2.00%  '
Set-TextValue $ws 'D7' '59.08'
Set-TextValue $ws 'E7' '  +9.32%  '
Set-TextValue $ws 'E8' '  -0.01%  '
Set-TextValue $ws 'E9' '  +4.06%  '
Set-TextValue $ws 'D10' '61.44'
Set-TextValue $ws 'E10' '  -0.04%  '
Set-TextValue $ws 'D11' '0.0800'
Set-TextValue $ws 'E11' '  +7.52%  '
Set-TextValue $ws 'D12' '0.108'
Set-TextValue $ws 'E12' '  +2.60%  '
Set-TextValue $ws 'D13' '16.25'
Set-TextValue $ws 'E13' '  +5.99%  '
Set-TextValue $ws 'D14' '2.373.72'
Set-TextValue $ws 'E14' '  -1.28%  '
Set-TextValue $ws 'D15' '0.821'
Set-TextValue $ws 'E15' '  -2.23%  '
Set-TextValue $ws 'D16' '5.53'
Set-TextValue $ws 'E16' '  +6.89%  '
Set-TextValue $ws 'D17' '2.073.81'
Set-TextValue $ws 'E17' '  -1.55%  '
Set-TextValue $ws 'D18' '37.047.74'
Set-TextValue $ws 'E18' '  -0.45%  '
Set-TextValue $ws 'D19' '15.76'
Set-TextValue $ws 'E19' '  +8.25%  '
Set-TextValue $ws 'B20' 'Litecoin'
Set-TextValue $ws 'C20' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws 'D20' '74.74'
Set-TextValue $ws 'E20' '  +3.06%  '
Set-TextValue $ws 'B21' 'ShibaInu'
Set-TextValue $ws 'C21' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws 'D21' '0.0₃0932'
Set-TextValue $ws 'E21' '  +11.26%  '
Set-TextValue $ws 'E22' '  +4.35%  '
Set-TextValue $ws 'D23' '240.36'
Set-TextValue $ws 'E23' '  -0.60%  '
Set-TextValue $ws 'E24' '  +0.04%  '
Set-TextValue $ws 'D25' '2.40'
Set-TextValue $ws 'E25' '  -2.95%  '
Set-TextValue $ws 'D26' '2.27'
Set-TextValue $ws 'E26' '  +13.31%  '
Set-TextValue $ws 'D27' '169.68'
Set-TextValue $ws 'E27' '  -0.60%  '
Set-TextValue $ws 'D28' '9.37'
Set-TextValue $ws 'E28' '  +1.29%  '
Set-TextValue $ws 'D29' '20.29'
Set-TextValue $ws 'E29' '  -1.90%  '
Set-TextValue $ws 'D30' '0.126'
Set-TextValue $ws 'E30' '  +2.64%  '
Set-TextValue $ws 'D31' '1.15'
Set-TextValue $ws 'E31' '  +6.16%  '
Set-TextValue $ws 'D32' '4.82'
Set-TextValue $ws 'E32' '  +6.41%  '
Set-TextValue $ws 'D33' '0.0635'
Set-TextValue $ws 'E33' '  +3.58%  '
Set-TextValue $ws 'D34' '4.49'
Set-TextValue $ws 'E34' '  +9.07%  '
Set-TextValue $ws 'D35' '0.0906'
Set-TextValue $ws 'E35' '  +0.17%  '
Set-TextValue $ws 'E36' '  +0.00%  '
Set-TextValue $ws 'D37' '2.31'
Set-TextValue $ws 'E37' '  -0.11%  '
Set-TextValue $ws 'E38' '  +29.10%  '
Set-TextValue $ws 'D39' '1.78'
Set-TextValue $ws 'E39' '  -3.91%  '
Set-TextValue $ws 'D41' '0.0227'
Set-TextValue $ws 'E41' '  +0.39%  '
Set-TextValue $ws 'D42' '17.73'
Set-TextValue $ws 'E42' '  -3.18%  '
Set-TextValue $ws 'E43' '  -0.40%  '
Set-TextValue $ws 'E44' '  +0.01%  '
Set-TextValue $ws 'D45' '4.43'
Set-TextValue $ws 'E45' '  +15.74%  '
Set-TextValue $ws 'E46' '  +1.50%  '
Set-TextValue $ws 'D47' '4.47'
Set-TextValue $ws 'E47' '  +11.30%  '
Set-TextValue $ws 'D48' '2.47'
Set-TextValue $ws 'E48' '  +8.36%  '
Set-TextValue $ws 'D49' '1.303.35'
Set-TextValue $ws 'E49' '  -1.46%  '
Set-TextValue $ws 'E50' '  -0.98%  '
Set-TextValue $ws 'D51' '6.95'
Set-TextValue $ws 'E51' '  -1.55%  '
